$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.608.91'
$ws.Range("E2").Value = '  -0.78%  '
$ws.Range("D3").Value = '3.482.74'
$ws.Range("E3").Value = '  -1.50%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.31'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.86%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.610'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.00%  '
$ws.Range("D9").Value = '3.482.80'
$ws.Range("E9").Value = '  -1.46%  '
$ws.Range("E10").Value = '  -1.85%  '
$ws.Range("E11").Value = '  -2.83%  '
$ws.Range("E12").Value = '  -3.56%  '
$ws.Range("D13").Value = '4.088.52'
$ws.Range("E13").Value = '  -1.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.39'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.16%  '
$ws.Range("E15").Value = '  -2.59%  '
$ws.Range("D16").Value = '67.594.62'
$ws.Range("E16").Value = '  -0.72%  '
$ws.Range("E17").Value = '  -2.93%  '
$ws.Range("D18").Value = '3.480.18'
$ws.Range("E18").Value = '  -1.34%  '
$ws.Range("E19").Value = '  -4.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.07'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '389.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.13%  '
$ws.Range("E22").Value = '  -1.82%  '
$ws.Range("E23").Value = '  +1.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.29%  '
$ws.Range("E26").Value = '  -2.20%  '
$ws.Range("E27").Value = '  -1.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.12'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.40%  '
$ws.Range("E29").Value = '  -1.29%  '
$ws.Range("E30").Value = '  +0.39%  '
$ws.Range("E31").Value = '  -4.92%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '24.67'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.05'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.44%  '
$ws.Range("E34").Value = '  -5.31%  '
$ws.Range("E35").Value = '  -3.77%  '
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("E37").Value = '  -4.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '160.58'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.889'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '27.92'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.77%  '
$ws.Range("E41").Value = '  -4.54%  '
$ws.Range("E42").Value = '  -4.40%  '
$ws.Range("E43").Value = '  -6.27%  '
$ws.Range("E44").Value = '  -4.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0712'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.04%  '
$ws.Range("D46").Value = '2.725.39'
$ws.Range("E46").Value = '  -6.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.98'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.19%  '
$ws.Range("E48").Value = '  -2.53%  '
$ws.Range("E49").Value = '  -3.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '330.82'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.11%  '
$ws.Range("E51").Value = '  -4.01%  '
